# Scheduled data refresh: update pricing/profit figures on each Leve sheet
# (currentAveragePrice*, Leve price/profit columns) with latest market data.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
# ALC row 33
$ws1.Range("H33").Value = 209.88235
$ws1.Range("I33").Value = 191.75
$ws1.Range("J33").Value = 500
$ws1.Range("K33").Value = 191.75
$ws1.Range("L33").Value = 500
$ws1.Range("M33").Value = 37.25
$ws1.Range("N33").Value = -958

# ALC row 53
$ws1.Range("H53").Value = 244
$ws1.Range("J53").Value = 200
$ws1.Range("L53").Value = 200
$ws1.Range("N53").Value = -1474

# ALC row 129
$ws1.Range("H129").Value = 295344.12
$ws1.Range("J129").Value = 323880.66
$ws1.Range("L129").Value = 971641.98
$ws1.Range("N129").Value = -981641.98

# ALC row 137
$ws1.Range("H137").Value = 2249.111
$ws1.Range("I137").Value = 2130
$ws1.Range("K137").Value = 6390
$ws1.Range("M137").Value = -3840

# ALC row 141
$ws1.Range("H141").Value = 2929.3157
$ws1.Range("I141").Value = 2350.4666
$ws1.Range("K141").Value = 7051.399800000001
$ws1.Range("M141").Value = -1871.399800000001

$ws2 = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws2.Range("H32").Value = 6671.9385
$ws2.Range("I32").Value = 4395.75
$ws2.Range("K32").Value = 4395.75
$ws2.Range("M32").Value = -4108.75

# ARM row 61
$ws2.Range("H61").Value = 1376.3478
$ws2.Range("I61").Value = 1375.1111
$ws2.Range("J61").Value = 1380.8
$ws2.Range("K61").Value = 1375.1111
$ws2.Range("L61").Value = 1380.8
$ws2.Range("M61").Value = -1163.1111
$ws2.Range("N61").Value = -1804.8

# ARM row 97
$ws2.Range("H97").Value = 90910200
$ws2.Range("I97").Value = 926.25
$ws2.Range("K97").Value = 926.25
$ws2.Range("M97").Value = -430.25

# ARM row 102
$ws2.Range("H102").Value = 2208.3333
$ws2.Range("J102").Value = 2765.8333
$ws2.Range("L102").Value = 2765.8333
$ws2.Range("N102").Value = -6009.8333

# ARM row 132
$ws2.Range("H132").Value = 14572.257
$ws2.Range("I132").Value = 1525.5483
$ws2.Range("K132").Value = 4576.644899999999
$ws2.Range("M132").Value = -2046.644899999999

# ARM row 136
$ws2.Range("H136").Value = 1376.3478
$ws2.Range("I136").Value = 1375.1111
$ws2.Range("J136").Value = 1380.8
$ws2.Range("K136").Value = 4125.3333
$ws2.Range("L136").Value = 4142.4
$ws2.Range("M136").Value = -1575.3333
$ws2.Range("N136").Value = -9242.4

$ws3 = $wb.Worksheets.Item("BSM")
# BSM row 64
$ws3.Range("H64").Value = 896.0625
$ws3.Range("I64").Value = 1109
$ws3.Range("J64").Value = 541.1667
$ws3.Range("K64").Value = 1109
$ws3.Range("L64").Value = 541.1667
$ws3.Range("M64").Value = -884
$ws3.Range("N64").Value = -991.1667

# BSM row 67
$ws3.Range("H67").Value = 896.0625
$ws3.Range("I67").Value = 1109
$ws3.Range("J67").Value = 541.1667
$ws3.Range("K67").Value = 1109
$ws3.Range("L67").Value = 541.1667
$ws3.Range("M67").Value = -329
$ws3.Range("N67").Value = -2101.1667

# BSM row 86
$ws3.Range("H86").Value = 1622.8223
$ws3.Range("I86").Value = 1395.1154
$ws3.Range("J86").Value = 1934.421
$ws3.Range("K86").Value = 1395.1154
$ws3.Range("L86").Value = 1934.421
$ws3.Range("M86").Value = -272.1153999999999
$ws3.Range("N86").Value = -4180.421

# BSM row 89
$ws3.Range("H89").Value = 1622.8223
$ws3.Range("I89").Value = 1395.1154
$ws3.Range("J89").Value = 1934.421
$ws3.Range("K89").Value = 6975.576999999999
$ws3.Range("L89").Value = 9672.105
$ws3.Range("M89").Value = -1359.576999999999
$ws3.Range("N89").Value = -20904.105

# BSM row 105
$ws3.Range("H105").Value = 3127493
$ws3.Range("I105").Value = 2954.4443
$ws3.Range("J105").Value = 7144757
$ws3.Range("K105").Value = 2954.4443
$ws3.Range("L105").Value = 7144757
$ws3.Range("M105").Value = -1207.4443
$ws3.Range("N105").Value = -7148251

$ws4 = $wb.Worksheets.Item("CRP")
# CRP row 31
$ws4.Range("H31").Value = 3714.6562
$ws4.Range("I31").Value = 2807.8572
$ws4.Range("J31").Value = 4419.9443
$ws4.Range("K31").Value = 2807.8572
$ws4.Range("L31").Value = 4419.9443
$ws4.Range("M31").Value = -2512.8572
$ws4.Range("N31").Value = -5009.9443

# CRP row 34
$ws4.Range("H34").Value = 3714.6562
$ws4.Range("I34").Value = 2807.8572
$ws4.Range("J34").Value = 4419.9443
$ws4.Range("K34").Value = 2807.8572
$ws4.Range("L34").Value = 4419.9443
$ws4.Range("M34").Value = -2605.8572
$ws4.Range("N34").Value = -4823.9443

# CRP row 86
$ws4.Range("H86").Value = 25827.666
$ws4.Range("I86").Value = 3966.6667
$ws4.Range("K86").Value = 3966.6667
$ws4.Range("M86").Value = -2843.6667

# CRP row 89
$ws4.Range("H89").Value = 25827.666
$ws4.Range("I89").Value = 3966.6667
$ws4.Range("K89").Value = 19833.3335
$ws4.Range("M89").Value = -14217.3335

# CRP row 132
$ws4.Range("H132").Value = 3007.08
$ws4.Range("I132").Value = 2240.75
$ws4.Range("K132").Value = 6722.25
$ws4.Range("M132").Value = -4192.25

# CRP row 134
$ws4.Range("H134").Value = 1366.6666
$ws4.Range("I134").Value = 1290.909
$ws4.Range("J134").Value = 1450
$ws4.Range("K134").Value = 3872.727
$ws4.Range("L134").Value = 4350
$ws4.Range("M134").Value = -1337.727
$ws4.Range("N134").Value = -9420

$ws5 = $wb.Worksheets.Item("CUL")
# CUL row 26
$ws5.Range("H26").Value = 529
$ws5.Range("J26").Value = 551
$ws5.Range("L26").Value = 1653
$ws5.Range("N26").Value = -2229

# CUL row 46
$ws5.Range("H46").Value = 1000
$ws5.Range("J46").Value = 1000
$ws5.Range("L46").Value = 3000
$ws5.Range("N46").Value = -3182

# CUL row 58
$ws5.Range("H58").Value = 3470.6667
$ws5.Range("I58").Value = 1000
$ws5.Range("J58").Value = 3964.8
$ws5.Range("K58").Value = 3000
$ws5.Range("L58").Value = 11894.4
$ws5.Range("M58").Value = -2872
$ws5.Range("N58").Value = -12150.4

# CUL row 69
$ws5.Range("H69").Value = 2500
$ws5.Range("I69").Value = 0
$ws5.Range("K69").Value = 0
$ws5.Range("M69").ClearContents()

# CUL row 72
$ws5.Range("H72").Value = 2500
$ws5.Range("I72").Value = 0
$ws5.Range("K72").Value = 0
$ws5.Range("M72").ClearContents()

# CUL row 103
$ws5.Range("H103").Value = 2732.8
$ws5.Range("I103").Value = 639
$ws5.Range("J103").Value = 4826.6
$ws5.Range("K103").Value = 1917
$ws5.Range("L103").Value = 14479.8
$ws5.Range("M103").Value = -1038
$ws5.Range("N103").Value = -16237.8

# CUL row 117
$ws5.Range("H117").Value = 1538.6666
$ws5.Range("I117").Value = 1300
$ws5.Range("J117").Value = 1777.3334
$ws5.Range("K117").Value = 3900
$ws5.Range("L117").Value = 5332.0002
$ws5.Range("M117").Value = -458
$ws5.Range("N117").Value = -12216.0002

# CUL row 121
$ws5.Range("H121").Value = 1050.2
$ws5.Range("J121").Value = 1118.1818
$ws5.Range("L121").Value = 3354.5454
$ws5.Range("N121").Value = -5974.5454

# CUL row 131
$ws5.Range("H131").Value = 719.88043
$ws5.Range("J131").Value = 724.43335
$ws5.Range("L131").Value = 2173.30005
$ws5.Range("N131").Value = -12253.30005

# CUL row 136
$ws5.Range("H136").Value = 4267.5
$ws5.Range("I136").Value = 1000
$ws5.Range("J136").Value = 4734.2856
$ws5.Range("K136").Value = 3000
$ws5.Range("L136").Value = 14202.8568
$ws5.Range("M136").Value = 2100
$ws5.Range("N136").Value = -24402.8568

# CUL row 138
$ws5.Range("H138").Value = 1999.2142
$ws5.Range("I138").Value = 1714.4546
$ws5.Range("J138").Value = 3043.3333
$ws5.Range("K138").Value = 5143.3638
$ws5.Range("L138").Value = 9129.999899999999
$ws5.Range("M138").Value = -3.363800000000083
$ws5.Range("N138").Value = -19409.9999

$ws6 = $wb.Worksheets.Item("GSM")
# GSM row 70
$ws6.Range("H70").Value = 3686000.5
$ws6.Range("I70").Value = 4274.875
$ws6.Range("J70").Value = 6958645.5
$ws6.Range("K70").Value = 4274.875
$ws6.Range("L70").Value = 6958645.5
$ws6.Range("M70").Value = -4004.875
$ws6.Range("N70").Value = -6959185.5

# GSM row 73
$ws6.Range("H73").Value = 3686000.5
$ws6.Range("I73").Value = 4274.875
$ws6.Range("J73").Value = 6958645.5
$ws6.Range("K73").Value = 4274.875
$ws6.Range("L73").Value = 6958645.5
$ws6.Range("M73").Value = -3338.875
$ws6.Range("N73").Value = -6960517.5

# GSM row 126
$ws6.Range("H126").Value = 5398.129
$ws6.Range("I126").Value = 4459.5293
$ws6.Range("J126").Value = 6537.857
$ws6.Range("K126").Value = 13378.5879
$ws6.Range("L126").Value = 19613.571
$ws6.Range("M126").Value = -10908.5879
$ws6.Range("N126").Value = -24553.571

# GSM row 132
$ws6.Range("H132").Value = 28873.736
$ws6.Range("I132").Value = 2216.2
$ws6.Range("K132").Value = 6648.599999999999
$ws6.Range("M132").Value = -4118.599999999999

$ws7 = $wb.Worksheets.Item("LTW")
# LTW row 68
$ws7.Range("H68").Value = 1975
$ws7.Range("I68").Value = 1900
$ws7.Range("J68").Value = 2200
$ws7.Range("K68").Value = 1900
$ws7.Range("L68").Value = 2200
$ws7.Range("M68").Value = -1151
$ws7.Range("N68").Value = -3698

# LTW row 71
$ws7.Range("H71").Value = 1975
$ws7.Range("I71").Value = 1900
$ws7.Range("J71").Value = 2200
$ws7.Range("K71").Value = 9500
$ws7.Range("L71").Value = 11000
$ws7.Range("M71").Value = -5756
$ws7.Range("N71").Value = -18488

# LTW row 132
$ws7.Range("H132").Value = 2113.8262
$ws7.Range("I132").Value = 1680
$ws7.Range("J132").Value = 2345.2
$ws7.Range("K132").Value = 5040
$ws7.Range("L132").Value = 7035.599999999999
$ws7.Range("M132").Value = -2510
$ws7.Range("N132").Value = -12095.6

$ws8 = $wb.Worksheets.Item("WVR")
# WVR row 64
$ws8.Range("H64").Value = 8972
$ws8.Range("I64").Value = 8972
$ws8.Range("K64").Value = 8972
$ws8.Range("M64").Value = -8724

# WVR row 67
$ws8.Range("H67").Value = 8972
$ws8.Range("I67").Value = 8972
$ws8.Range("K67").Value = 8972
$ws8.Range("M67").Value = -8114

# WVR row 132
$ws8.Range("H132").Value = 1087.2307
$ws8.Range("I132").Value = 642.0645
$ws8.Range("K132").Value = 1926.1935
$ws8.Range("M132").Value = 603.8065000000001

# WVR row 136
$ws8.Range("H136").Value = 21741398
$ws8.Range("I136").Value = 27028058
$ws8.Range("J136").Value = 7348.8887
$ws8.Range("K136").Value = 81084174
$ws8.Range("L136").Value = 22046.6661
$ws8.Range("M136").Value = -81081624
$ws8.Range("N136").Value = -27146.6661
